$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = '26.386.95'
$ws.Cells.Item(2, 5).Value = '  -1.92%  '

$ws.Cells.Item(3, 4).Value = '1.837.05'
$ws.Cells.Item(3, 5).Value = '  -2.21%  '

$ws.Cells.Item(4, 4).Value = "'1.000"
$ws.Cells.Item(4, 4).Style = 'Normal'
$ws.Cells.Item(4, 5).Value = '  +0.00%  '

$ws.Cells.Item(5, 4).Value = "'259.82"
$ws.Cells.Item(5, 4).Style = 'Normal'
$ws.Cells.Item(5, 5).Value = '  -6.85%  '

$ws.Cells.Item(6, 5).Value = '  +0.00%  '

$ws.Cells.Item(7, 4).Value = "'0.5220"
$ws.Cells.Item(7, 4).Style = 'Normal'

$ws.Cells.Item(8, 4).Value = "'0.3226"
$ws.Cells.Item(8, 4).Style = 'Normal'
$ws.Cells.Item(8, 5).Value = '  -6.64%  '

$ws.Cells.Item(9, 4).Value = "'0.06755"
$ws.Cells.Item(9, 4).Style = 'Normal'
$ws.Cells.Item(9, 5).Value = '  -2.94%  '

$ws.Cells.Item(10, 4).Value = "'18.56"
$ws.Cells.Item(10, 4).Style = 'Normal'
$ws.Cells.Item(10, 5).Value = '  -7.87%  '

$ws.Cells.Item(11, 4).Value = "'0.7633"
$ws.Cells.Item(11, 4).Style = 'Normal'
$ws.Cells.Item(11, 5).Value = '  -5.35%  '

$ws.Cells.Item(12, 4).Value = "'0.07677"
$ws.Cells.Item(12, 4).Style = 'Normal'
$ws.Cells.Item(12, 5).Value = '  -2.36%  '

$ws.Cells.Item(13, 4).Value = '1.864.74'
$ws.Cells.Item(13, 5).Value = '  -0.61%  '

$ws.Cells.Item(14, 4).Value = "'88.50"
$ws.Cells.Item(14, 4).Style = 'Normal'
$ws.Cells.Item(14, 5).Value = '  -1.98%  '

$ws.Cells.Item(15, 4).Value = "'5.013"
$ws.Cells.Item(15, 4).Style = 'Normal'
$ws.Cells.Item(15, 5).Value = '  -3.24%  '

$ws.Cells.Item(16, 5).Value = '  +0.07%  '

$ws.Cells.Item(17, 4).Value = "'13.94"
$ws.Cells.Item(17, 4).Style = 'Normal'
$ws.Cells.Item(17, 5).Value = '  -4.75%  '

$ws.Cells.Item(18, 5).Value = '  +0.00%  '

$ws.Cells.Item(19, 4).Value = "'0.000007915"
$ws.Cells.Item(19, 4).Style = 'Normal'
$ws.Cells.Item(19, 5).Value = '  -2.39%  '

$ws.Cells.Item(20, 4).Value = '26.435.59'
$ws.Cells.Item(20, 5).Value = '  -1.93%  '

$ws.Cells.Item(21, 4).Value = '2.074.37'
$ws.Cells.Item(21, 5).Value = '  -1.80%  '

$ws.Cells.Item(22, 4).Value = "'4.559"
$ws.Cells.Item(22, 4).Style = 'Normal'
$ws.Cells.Item(22, 5).Value = '  -4.15%  '

$ws.Cells.Item(23, 4).Value = "'9.432"
$ws.Cells.Item(23, 4).Style = 'Normal'
$ws.Cells.Item(23, 5).Value = '  -6.27%  '

$ws.Cells.Item(24, 4).Value = "'5.926"
$ws.Cells.Item(24, 4).Style = 'Normal'
$ws.Cells.Item(24, 5).Value = '  -4.33%  '

$ws.Cells.Item(25, 4).Value = "'144.86"
$ws.Cells.Item(25, 4).Style = 'Normal'
$ws.Cells.Item(25, 5).Value = '  -0.83%  '

$ws.Cells.Item(26, 4).Value = "'2.223"
$ws.Cells.Item(26, 4).Style = 'Normal'
$ws.Cells.Item(26, 5).Value = '  -5.64%  '

$ws.Cells.Item(27, 4).Value = "'1.654"
$ws.Cells.Item(27, 4).Style = 'Normal'
$ws.Cells.Item(27, 5).Value = '  -0.53%  '

$ws.Cells.Item(28, 4).Value = "'16.95"
$ws.Cells.Item(28, 4).Style = 'Normal'
$ws.Cells.Item(28, 5).Value = '  -2.40%  '

$ws.Cells.Item(29, 4).Value = "'111.36"
$ws.Cells.Item(29, 4).Style = 'Normal'
$ws.Cells.Item(29, 5).Value = '  -2.37%  '

$ws.Cells.Item(30, 4).Value = "'4.168"
$ws.Cells.Item(30, 4).Style = 'Normal'
$ws.Cells.Item(30, 5).Value = '  -4.77%  '

$ws.Cells.Item(31, 4).Value = "'4.136"
$ws.Cells.Item(31, 4).Style = 'Normal'
$ws.Cells.Item(31, 5).Value = '  -4.64%  '

$ws.Cells.Item(32, 4).Value = "'0.08745"
$ws.Cells.Item(32, 4).Style = 'Normal'
$ws.Cells.Item(32, 5).Value = '  -1.84%  '

$ws.Cells.Item(33, 4).Value = "'0.04816"
$ws.Cells.Item(33, 4).Style = 'Normal'
$ws.Cells.Item(33, 5).Value = '  -2.59%  '

$ws.Cells.Item(34, 4).Value = "'1.120"
$ws.Cells.Item(34, 4).Style = 'Normal'
$ws.Cells.Item(34, 5).Value = '  -5.02%  '

$ws.Cells.Item(35, 4).Value = "'2.843"
$ws.Cells.Item(35, 4).Style = 'Normal'
$ws.Cells.Item(35, 5).Value = '  -1.43%  '

$ws.Cells.Item(36, 4).Value = "'0.7009"
$ws.Cells.Item(36, 4).Style = 'Normal'
$ws.Cells.Item(36, 5).Value = '  -4.92%  '

$ws.Cells.Item(37, 4).Value = "'3.056"
$ws.Cells.Item(37, 4).Style = 'Normal'
$ws.Cells.Item(37, 5).Value = '  -6.82%  '

$ws.Cells.Item(38, 4).Value = "'0.01761"
$ws.Cells.Item(38, 4).Style = 'Normal'
$ws.Cells.Item(38, 5).Value = '  -4.89%  '

$ws.Cells.Item(39, 4).Value = "'2.190"
$ws.Cells.Item(39, 4).Style = 'Normal'
$ws.Cells.Item(39, 5).Value = '  -8.27%  '

$ws.Cells.Item(40, 4).Value = "'0.4827"
$ws.Cells.Item(40, 4).Style = 'Normal'
$ws.Cells.Item(40, 5).Value = '  -6.44%  '

$ws.Cells.Item(41, 4).Value = "'111.24"
$ws.Cells.Item(41, 4).Style = 'Normal'
$ws.Cells.Item(41, 5).Value = '  -4.66%  '

$ws.Cells.Item(42, 4).Value = "'0.8850"
$ws.Cells.Item(42, 4).Style = 'Normal'
$ws.Cells.Item(42, 5).Value = '  -7.49%  '

$ws.Cells.Item(43, 4).Value = "'6.090"
$ws.Cells.Item(43, 4).Style = 'Normal'
$ws.Cells.Item(43, 5).Value = '  -1.99%  '

$ws.Cells.Item(44, 4).Value = "'0.9992"
$ws.Cells.Item(44, 4).Style = 'Normal'
$ws.Cells.Item(44, 5).Value = '  +0.01%  '

$ws.Cells.Item(45, 4).Value = "'7.620"
$ws.Cells.Item(45, 4).Style = 'Normal'
$ws.Cells.Item(45, 5).Value = '  -6.43%  '

$ws.Cells.Item(46, 2).Value = 'Cronos'
$ws.Cells.Item(46, 3).Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Cells.Item(46, 4).Value = "'0.05851"
$ws.Cells.Item(46, 4).Style = 'Normal'
$ws.Cells.Item(46, 5).Value = '  -1.83%  '

$ws.Cells.Item(47, 2).Value = 'Decentraland'
$ws.Cells.Item(47, 3).Value = 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
$ws.Cells.Item(47, 4).Value = "'0.4107"
$ws.Cells.Item(47, 4).Style = 'Normal'
$ws.Cells.Item(47, 5).Value = '  -8.77%  '

$ws.Cells.Item(48, 4).Value = "'9.000"
$ws.Cells.Item(48, 4).Style = 'Normal'
$ws.Cells.Item(48, 5).Value = '  -4.16%  '

$ws.Cells.Item(49, 4).Value = "'34.71"
$ws.Cells.Item(49, 4).Style = 'Normal'
$ws.Cells.Item(49, 5).Value = '  -4.23%  '

$ws.Cells.Item(50, 5).Value = '  -9.50%  '

$ws.Cells.Item(51, 4).Value = "'0.8796"
$ws.Cells.Item(51, 4).Style = 'Normal'
$ws.Cells.Item(51, 5).Value = '  -0.77%  '
